$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1 and G1, matching the style of the existing header row (A1:E1):
# bold font, centered horizontal alignment.
$ws.Range("F1").Value = "conf.low"
$ws.Range("G1").Value = "conf.high"
$ws.Range("F1:G1").Font.Bold = $true
$ws.Range("F1:G1").HorizontalAlignment = -4108

# Fill in the conf.low / conf.high values for rows 2-12
$confLow = @(
    0.5404943604321496,
    -0.1734842019197602,
    -0.3096297438804556,
    -0.1657064096745759,
    0.05961719515769731,
    -0.05771013484998045,
    -0.00169392465957928,
    -0.03018175431607399,
    -0.05047951731765825,
    -0.02303032095436251,
    -0.05814971966688868
)

$confHigh = @(
    0.6930403085996321,
    -0.09139143722564977,
    -0.2229285242962351,
    -0.07968324298859122,
    0.1449152733077592,
    0.03286243204005899,
    0.0009545036779651839,
    0.04122469103777791,
    0.0280754315198984,
    0.0555499096555235,
    0.03076935272063646
)

for ($i = 0; $i -lt $confLow.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $confLow[$i]
    $ws.Cells.Item($row, 7).Value = $confHigh[$i]
}
